$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32-47 (COS port entries 30:0 .. 37:1) currently repeat the same
# "Dispositivo DAAS" / "Puerto DAAS" / "Unnamed: 5" values
# (MEDE-CABA-H-03-DAAS / xe-0/0/* / PUERTOLIBRE). These were renumbered away
# and the three columns (E:G) for those rows should end up blank, matching
# the pattern already used from row 48 onward.
#
# A plain ClearContents()/Value = "" removes the cell entirely (making it a
# true empty/blank cell), but the target state keeps the cells present as
# empty text cells. Writing the Excel "text prefix" apostrophe forces the
# cell to be stored as an empty string (text) instead of being dropped, and
# re-applying the "Normal" style afterwards strips the quote-prefix
# formatting Excel would otherwise tag the cell with.
$target = $ws.Range("E32:G47")
$target.Value = "'"
$target.Style = "Normal"
